# Updated benchmarks xls; wrapper to produce distr benchmarks hangs; new paper version
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2: Guacamol value 0.979 -> 0.999969, recolor from yellow to green (reuses
# the existing "green" style/fill instead of adding a new one).
$ws.Range("C2").Value = 0.999969
$ws.Range("C2").Interior.Color = 5296274   # RGB(80,208,146) == FF92D050 -> green fill

# C7: Guacamol value 0.75 -> 0.961, recolor from orange to yellow (reuses the
# existing "yellow" style/fill).
$ws.Range("C7").Value = 0.961
$ws.Range("C7").Interior.Color = 65535     # RGB(255,255,0) == FFFFFF00 -> yellow fill

# C8: Guacamol value 0.99722 -> 0.99741 (style/fill unchanged - stays yellow).
$ws.Range("C8").Value = 0.99741000000000002

# Move the active selection from C17 to B20, as last left by the author.
[void]$ws.Range("B20").Select()
